$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value. Values are written as text (matching the
# workbooks inline-string storage for these columns), using column C (a
# plain text column with the sheets default style) as the style template
# so no stray number-format styling is left behind on the edited cells.
function Set-TextCell($row, $col, $value) {
    $cell = $ws.Range("$col$row")
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $ws.Range("C$row").Style
}

# Row 2
Set-TextCell 2 "D" "319.19"
Set-TextCell 2 "E" "4.31%"

# Row 3
Set-TextCell 3 "D" "39.58"
Set-TextCell 3 "E" "3.27%"

# Row 4
Set-TextCell 4 "E" "0.75%"

# Row 5
Set-TextCell 5 "D" "0.08194"
Set-TextCell 5 "E" "1.62%"

# Row 6
Set-TextCell 6 "D" "2.018"
Set-TextCell 6 "E" "3.73%"

# Row 7
Set-TextCell 7 "D" "8.278"
Set-TextCell 7 "E" "4.34%"

# Row 8
Set-TextCell 8 "D" "4.275"
Set-TextCell 8 "E" "2.16%"

# Row 9
Set-TextCell 9 "D" "0.9340"
Set-TextCell 9 "E" "0.57%"

# Row 10
Set-TextCell 10 "D" "0.1405"
Set-TextCell 10 "E" "-4.23%"

# Row 11
Set-TextCell 11 "D" "0.1997"
Set-TextCell 11 "E" "3.49%"

# Row 12
Set-TextCell 12 "D" "0.09095"
Set-TextCell 12 "E" "1.00%"

# Row 13
Set-TextCell 13 "D" "0.03582"
Set-TextCell 13 "E" "2.52%"

# Row 14
Set-TextCell 14 "D" "0.09810"
Set-TextCell 14 "E" "0.24%"

# Row 15
Set-TextCell 15 "D" "0.001394"
Set-TextCell 15 "E" "0.27%"

# Row 16
Set-TextCell 16 "D" "0.006008"
Set-TextCell 16 "E" "2.15%"

# Row 17
Set-TextCell 17 "D" "3.661"
Set-TextCell 17 "E" "-1.62%"

# Row 18
Set-TextCell 18 "E" "-5.42%"

# Row 19
Set-TextCell 19 "E" "-0.06%"

# Row 20
Set-TextCell 20 "E" "-2.13%"

# Row 21
Set-TextCell 21 "D" "4.901"
Set-TextCell 21 "E" "1.95%"

# Row 22
Set-TextCell 22 "E" "1.94%"

# Row 23
Set-TextCell 23 "D" "0.04331"
Set-TextCell 23 "E" "-0.81%"

# Row 24
Set-TextCell 24 "E" "-0.57%"

# Row 25
Set-TextCell 25 "D" "0.004785"
Set-TextCell 25 "E" "12.27%"

# Row 26
Set-TextCell 26 "E" "0.05%"

# Row 27
Set-TextCell 27 "D" "0.0004001"
Set-TextCell 27 "E" "-10.03%"

# Row 39
Set-TextCell 39 "E" "7.01%"

# Row 40
Set-TextCell 40 "D" "0.05261"
Set-TextCell 40 "E" "4.31%"

# Row 41
Set-TextCell 41 "D" "0.007526"
Set-TextCell 41 "E" "0.74%"

# Row 42
Set-TextCell 42 "D" "0.01013"
Set-TextCell 42 "E" "0.24%"

# Row 43
Set-TextCell 43 "D" "0.1381"
Set-TextCell 43 "E" "2.35%"

# Row 44
Set-TextCell 44 "E" "0.52%"

# Row 45
Set-TextCell 45 "D" "0.009878"
Set-TextCell 45 "E" "10.62%"

# Row 46
Set-TextCell 46 "D" "0.00006498"
Set-TextCell 46 "E" "5.16%"

# Row 47
Set-TextCell 47 "E" "0.09%"

# Row 48
Set-TextCell 48 "E" "-1.20%"

# Row 49
Set-TextCell 49 "E" "-24.96%"

# Row 50
Set-TextCell 50 "D" "0.00002101"
Set-TextCell 50 "E" "0.09%"

# Row 51
Set-TextCell 51 "D" "0.0002001"
Set-TextCell 51 "E" "0.09%"
